$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "classificacao_painel" - classification of the instrument in the Painel
$ws.Range("H1").Value = "classificacao_painel"
$ws.Range("H2").Value = "Ampliação do Bolsa Família"
$ws.Range("H3").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H4").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H5").Value = "Transferência para a conta de Desenvolvimento Energético"
$ws.Range("H6").Value = "Transferência para a conta de Desenvolvimento Energético"
$ws.Range("H7").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H8").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H9").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H10").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H11").Value = "Auxílio Financeiro aos Estados, Municípios e DF"
$ws.Range("H12").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H13").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H14").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H15").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H16").Value = "Benefício Emergencial de Manutenção do Emprego e Renda"
$ws.Range("H17").Value = "Concessão de Financiamento para pagamento de folha salarial"
$ws.Range("H18").Value = "Concessão de Financiamento para pagamento de folha salarial - Operacionalização financeira pela MP 943/2020"
$ws.Range("H19").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H20").Value = "Financiamento de Infraestrutura Turística"
$ws.Range("H21").Value = "Cotas dos Fundos Garantidores de Operações e de Crédito"
$ws.Range("H22").Value = "Cotas dos Fundos Garantidores de Operações e de Crédito"
$ws.Range("H23").Value = "Auxílio Emergencial a pessoas em situação de vulnerabilidade, devido à pandemia da Covid-19"
$ws.Range("H24").Value = "Cotas dos Fundos Garantidores de Operações e de Crédito"
$ws.Range("H25").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H26").Value = "Programa Emergencial de Acesso a Crédito - Maquininhas"
$ws.Range("H27").Value = "Cotas dos Fundos Garantidores de Operações e de Crédito"
$ws.Range("H28").Value = "Benefício Emergencial de Manutenção do Emprego e Renda"
$ws.Range("H29").Value = "Não está no Painel, pois trata apenas de regras. A operacionalização financeira foi feita pela MP nº 939/2020"
$ws.Range("H30").Value = "Auxílio Financeiro aos Estados, Municípios e DF"
$ws.Range("H31").Value = "Auxílio Financeiro aos Estados, Municípios e DF"
$ws.Range("H32").Value = "Auxílio Financeiro aos Estados, Municípios e DF"
$ws.Range("H33").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H34").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H35").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H36").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H37").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H38").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H39").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H40").Value = "Auxílio Financeiro aos Estados, Municípios e DF"
$ws.Range("H41").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H42").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H43").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H44").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H45").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H46").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H47").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H48").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H49").Value = "Aquisição de Vacinas"
$ws.Range("H50").Value = "Aquisição de Vacinas"
$ws.Range("H51").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H52").Value = "Aquisição de Vacinas"
$ws.Range("H53").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H54").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H55").Value = "Despesas Adicionais do Ministério da Saúde e demais ministérios"
$ws.Range("H56").Value = "Aquisição de Vacinas *reclassificado no Painel no dia 21 de junho de 2021"

# Widen the new column to fit its content
$ws.Columns.Item(8).ColumnWidth = 61.67

# Restore the view: no fixed top-left cell, adjusted zoom, selection moved to A2
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.Zoom = 110
